$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "cat" row (row 22) to the data dictionary ---

# Row height matches the other wrapped/merged description rows
$ws.Rows(22).RowHeight = 32.25

# Cell values
$ws.Range("A22").Value = "cat"
$ws.Range("B22").Value = "Categoria"
$ws.Range("C22").Value = "Variable de predicción"
$ws.Range("D22").Value = "0 Barato`n1 Caro"

# Formatting: re-use the existing formatting already used elsewhere in the
# sheet so no unrelated styles are disturbed.
# A22 -> plain/default style (same as most NOMBRE cells, e.g. A2)
$ws.Range("A2").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "cat"

# B22/C22 -> same style as the other ETIQUETA/PREGUNTA cells (e.g. row 3)
$ws.Range("B3").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C22").PasteSpecial(-4122)

# D22 -> same wrapped style used by the other POSIBLES RESPUESTAS cells (e.g. D6)
$ws.Range("D6").Copy()
$ws.Range("D22").PasteSpecial(-4122)

# Merge the answer cell across D:F like every other row
$ws.Range("D22:F22").Merge()

# E22/F22 -> plain/default style for the (now merged, blank) trailing cells
$ws.Range("A2").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").PasteSpecial(-4122)

# Update the selection to reflect where the cursor ended up after the edit
$ws.Range("D23").Select()
